$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Libraries")

# Add the two new library/suffix rows (99c / D1SUNED FENDL3.1d+EAF2007 and 93c / D1SUNED FENDL3.2b+DECAY2020)
$ws.Range("A9").Value = "99c"
$ws.Range("B9").Value = "D1SUNED (FENDL 3.1d+EAF2007)"
$ws.Range("A10").Value = "93c"
$ws.Range("B10").Value = "D1SUNED (FENDL 3.2b+DECAY2020)"

# Column B now holds longer text, widen it to fit
$ws.Columns.Item(2).ColumnWidth = 31.25

# Make Libraries the active sheet/tab with the same selected cell as the source workbook
$ws.Activate() | Out-Null
$ws.Range("C19").Select() | Out-Null
